$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the dSF column (F) values per repulled data / mean calculation
$ws.Range("F2").Value = -12
$ws.Range("F3").Value = -2
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 0
